$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap row 3 and row 4 contents (Brown coal <-> Black coal and their values)
$ws.Range("A3").Value = "Black coal"
$ws.Range("B3").Value = -4.283866833183724
$ws.Range("C3").Value = 12383.7065
$ws.Range("D3").Value = 9114.790000000001

$ws.Range("A4").Value = "Brown coal"
$ws.Range("B4").Value = -3.237219697661808
$ws.Range("C4").Value = 4597.9545
$ws.Range("D4").Value = 3651.935
